$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Text columns (Coin name / Link) - plain strings, no numeric coercion risk
$ws.Range("B9").Value = "One"
$ws.Range("C9").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
$ws.Range("B15").Value = "TigerCash"
$ws.Range("C15").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("B16").Value = "LEO"
$ws.Range("C16").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("B17").Value = "GateToken"
$ws.Range("C17").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("B18").Value = "BTSEToken"
$ws.Range("C18").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
$ws.Range("B19").Value = "BitpandaEcosystemToken"
$ws.Range("C19").Value = "https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best"
$ws.Range("B20").Value = "WazirX"
$ws.Range("C20").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"

# Numeric-looking text columns (Price / Volume%) - force text format first
# so Excel does not coerce these into Number values, matching the source
# workbook where these cells are stored as text.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "245.32"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "-0.84%"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "27.25"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "3.74%"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.105"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "0.33%"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.05686"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "1.52%"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "0.80%"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.8198"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "0.83%"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.8609"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "2.02%"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.0005990"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "0.33%"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "-0.55%"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.02853"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "1.48%"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.09391"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "-0.05%"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.001513"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "-0.19%"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.04075"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "-13.11%"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.006215"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "0.43%"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.511"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "-2.66%"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.010"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "-0.30%"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.317"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "12.72%"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.3164"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "1.37%"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.1333"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "-0.41%"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.03226"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "0.91%"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "-0.09%"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "3.563"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "-4.84%"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "1.77%"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.001217"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "-2.08%"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.004467"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "-2.35%"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "22.95%"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "-25.42%"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "1.57%"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.005930"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "72.90%"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "-21.86%"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.002300"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "-13.51%"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.009716"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "17.86%"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00005116"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "-5.28%"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "-30.35%"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "-2.89%"
